# Add a new calculated column "Less than 10 units" to Table1 on the
# SalesOrders sheet, flagging rows whose Units value is below 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SalesOrders")
$lo = $ws.ListObjects.Item("Table1")

# Grow the table by one column and name it.
$col = $lo.ListColumns.Add()
$headerCell = $lo.HeaderRowRange.Cells.Item(1, $col.Index)
$headerCell.Value = "Less than 10 units"

# Fill the new column with a calculated-column formula referencing Units.
$lo.ListColumns.Item($col.Index).DataBodyRange.Formula = "=Table1[[#This Row],[Units]]<10"

# Give the column a sensible display width.
$ws.Columns.Item($col.Index).AutoFit() | Out-Null

# Make SalesOrders the active sheet/selection, matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("H2").Select() | Out-Null
